# Update the workbook's build/version string everywhere it appears.
#
# Old version string: "mines - January 30 (built on February 02 2026 12.49.33 EST)"
# New version string: "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet updates ---

# A2: "Version: mines - January 30 (...)" -> "Version: mines - version 1.0.0 (...)"
$aboutSheet.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended citation string containing the version in quotes.
$aboutSheet.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for Zhangji Coal Mine, China, M0426, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet updates ---
# Column S ("build_version") rows 2 through 11 hold the same version string.

for ($row = 2; $row -le 11; $row++) {
    $dataSheet.Cells.Item($row, 19).Value = $newVersion
}
